# ALIGNED: ram request will have 64bit aligned address
#
# 1) Bump the cached "datetimeFigureOut" footer field from 2020/10/12 to
#    2020/10/13 everywhere it appears (slide master + every slide layout).
# 2) Fix the mislabeled "Class RAM {" header that sits under the CPU
#    pseudocode box on slide 1 so it reads "Class CPU {" (split into two
#    runs: "Class CPU " + "{", matching how the existing tail of the
#    original run is preserved).

$p = $ppt.ActivePresentation

function Update-DateShape {
    param($shape)

    if ($shape.HasTextFrame -ne -1) {
        return
    }

    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -eq "2020/10/12") {
        $tr.Text = "2020/10/13"
    }
}

# --- Slide master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# --- Every slide layout hanging off the master ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# --- Slide 1: relabel the CPU class box ---
# Both the CPU-side and RAM-side pseudocode boxes currently start with
# "Class RAM {" (copy/paste leftover) — only the left-hand one (smaller
# .Left, sitting under the "CPU" caption) is mislabeled and needs fixing.
$slide = $p.Slides.Item(1)
$target = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -ne -1) {
        continue
    }
    $tr = $shape.TextFrame.TextRange
    if ($tr.Text.IndexOf("Class RAM {") -eq 0) {
        if (($target -eq $null) -or ($shape.Left -lt $target.Left)) {
            $target = $shape
        }
    }
}

if ($target -ne $null) {
    $tr = $target.TextFrame.TextRange
    # "Class RAM {" -> "Class CPU {", splitting off the leading
    # "Class CPU " into its own run and leaving the trailing "{"
    # as the remainder of the original run.
    $head = $tr.Characters(1, 10)
    $head.Text = "Class CPU "
}
